$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "NAQUICHE SILVA MARIA LIZETH",
    "NAQUICHE MECHATO ENMA",
    "MANOSALVA RUIZ SANDRA KAROLINE",
    "CUBAS GARCIA ROSA ANITA",
    "CORAS QUISPE JORGE AMERICO",
    "BECERRA ASMAT CAROL STEFANY",
    "SAUCEDO CABRERA CARLOS ALEXANDER",
    "CASTREJON TELLO GRECIA",
    "PACHECO ALISON",
    "45752721",
    "44591379"
)

$totals = @(80, 75, 62, 62, 57, 50, 34, 20, 17, 1, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    if ($names[$i] -eq "45752721" -or $names[$i] -eq "44591379") {
        # Force a text cell (these look numeric but must remain text, as
        # in the source) then strip the formatting back to Normal so no
        # extra style survives on the cell.
        $cellA.NumberFormat = "@"
        $cellA.Value = $names[$i]
        $cellA.Style = "Normal"
    } else {
        $cellA.Value = $names[$i]
    }
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
